$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.341.99'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.21%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.488.18'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.86%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '597.23'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.33%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '176.93'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.41%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.01%  '

# Row 8
$ws.Range('E8').Value = '  -0.19%  '

# Row 9
$ws.Range('E9').Value = '  -2.30%  '

# Row 10
$ws.Range('E10').Value = '  -2.73%  '

# Row 11
$ws.Range('E11').Value = '  -2.07%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.088.05'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.11%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '31.61'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +10.20%  '

# Row 14
$ws.Range('E14').Value = '  -0.03%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '67.327.14'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.20%  '

# Row 16
$ws.Range('E16').Value = '  -2.87%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.480.71'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.18%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.25'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.40%  '

# Row 19
$ws.Range('E19').Value = '  +1.26%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '388.91'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.84%  '

# Row 21
$ws.Range('E21').Value = '  -0.61%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '72.92'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.67%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.05%  '

# Row 24
$ws.Range('E24').Value = '  -0.78%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.72'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.01%  '

# Row 26
$ws.Range('E26').Value = '  -0.57%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.27'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.31%  '

# Row 28
$ws.Range('E28').Value = '  -1.28%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.994'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.64%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.16'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.00%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.43'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.52%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.05'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.66%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '23.59'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.12%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '7.29'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.51%  '

# Row 35
$ws.Range('E35').Value = '  +0.06%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '163.67'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.05%  '

# Row 37
$ws.Range('E37').Value = '  +0.91%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.871'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.66%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.99'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.17%  '

# Row 40
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.64'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.59%  '

# Row 41
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '27.24'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.64%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '26.47'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.09%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.816.46'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.52%  '

# Row 44
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0723'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.10%  '

# Row 45
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.58'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.58%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '42.24'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.65%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '341.60'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.13%  '

# Row 48
$ws.Range('E48').Value = '  -3.35%  '

# Row 49
$ws.Range('E49').Value = '  -2.70%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '33.35'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.73%  '

# Row 51
$ws.Range('E51').Value = '  -2.39%  '

Write-Output "applied cryptos update"